$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - new record appended by the Streamlit automation
$ws.Range("A5").Value = "DF"
$ws.Range("B5").Value = "03AN313"

# SUB1/SUB2/SUB3 are present but blank on this row, same as the rows above it.
# Force a real (empty) text cell instead of letting an empty assignment
# remove the cell entirely.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = ""
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = ""
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = ""

$ws.Range("F5").Value = "BATERIA"

# MAQUINAS ("8334") must stay text (as in the source data), not become a number.
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "8334"

$ws.Range("H5").Value = "SICOOB - (78KLM10 13/11/25_24/03/2020) - DF"
$ws.Range("I5").Value = "13/11/25"
$ws.Range("J5").Value = "24/03/2020"
$ws.Range("K5").Value = "13/11/25"
$ws.Range("L5").Value = "DENTRO"

# Drop the temporary text format so the new cells fall back to the
# workbook's default (unstyled) look, matching the rest of the sheet.
$ws.Range("C5:E5").Style = "Normal"
$ws.Range("G5").Style = "Normal"
